# The commit removes the "Age" column (column C) from both the "Valence"
# and "Arousal" worksheets. Deleting the entire column shifts the
# subsequent columns (Fearful, Neutral, Happy) left by one and Excel
# automatically re-targets every formula/reference that pointed at the
# shifted ranges (AVERAGE/STDEV.S/SEM rows, shared-formula ranges, etc.).
# It also drops the now-unused "Age" shared string.

$wb = $excel.ActiveWorkbook

$wsValence = $wb.Worksheets.Item("Valence")
$wsArousal = $wb.Worksheets.Item("Arousal")

# Remove the Age column (column C) on both sheets.
$wsValence.Range("C1").EntireColumn.Delete()
$wsArousal.Range("C1").EntireColumn.Delete()

# Match the author's final selection / active-sheet state as closely as
# the object model allows: the workbook was left with "Arousal" as the
# active tab, cell K18 selected there, and cell M35 selected (but not
# active) back on "Valence".
$wsValence.Range("M35").Select()
$wsArousal.Activate()
$wsArousal.Range("K18").Select()
